# Update graphs and correlations
#
# - "  Nepal" (A119, leading spaces) -> "Nepal"
# - " Switzerland" (A166, leading space) -> "Switzerland"
# - "United States of America" (A181) -> "United States"
# - Append a new "World" row (row 191) with values 41.3 / 41.3
# - Update the view state (selection + scroll position) to match the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string slot assignment on save (matches the
# original author's edit order so the resulting sharedStrings.xml indices
# line up exactly with the target file).
$ws.Range("A181").Value = "United States"
$ws.Range("A166").Value = "Switzerland"
$ws.Range("A119").Value = "Nepal"

$ws.Range("A191").Value = "World"
$ws.Range("B191").Value = 41.3
$ws.Range("C191").Value = 41.3

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 161
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B192").Select()
